$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.292879
$ws.Range("H2").Value = 0.878637
$ws.Range("I2").Value = 0.007455471808396097
$ws.Range("J2").Value = 0.007455471808396097
$ws.Range("M2").Value = 0.1124023333333333
$ws.Range("N2").Value = 0.337207
$ws.Range("O2").Value = 0.03490487583665934
$ws.Range("P2").Value = 0.03490487583665934
$ws.Range("Q2").Value = 0.03292028298433333
$ws.Range("R2").Value = 0.296282546859
$ws.Range("S2").Value = 0.0002602323177757798
$ws.Range("T2").Value = 0.0002602323177757798

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.292879
$ws.Range("H3").Value = 0.878637
$ws.Range("I3").Value = 0.007455471808396097
$ws.Range("J3").Value = 0.007455471808396097
$ws.Range("O3").Value = 0.1007904982867776
$ws.Range("P3").Value = 0.1007904982867776
$ws.Range("Q3").Value = 0.09505983465633334
$ws.Range("R3").Value = 0.855538511907
$ws.Range("S3").Value = 0.0007514407185312651
$ws.Range("T3").Value = 0.0007514407185312652

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.292879
$ws.Range("H4").Value = 0.878637
$ws.Range("I4").Value = 0.007455471808396097
$ws.Range("J4").Value = 0.007455471808396097
$ws.Range("M4").Value = 2.783274666666667
$ws.Range("N4").Value = 8.349824
$ws.Range("O4").Value = 0.8643046258765631
$ws.Range("P4").Value = 0.8643046258765632
$ws.Range("Q4").Value = 0.8151627010986667
$ws.Range("R4").Value = 7.336464309888
$ws.Range("S4").Value = 0.006443798772089051
$ws.Range("T4").Value = 0.006443798772089052

# Row 5
$ws.Range("I5").Value = 0.6729737392616156
$ws.Range("J5").Value = 0.6729737392616155
$ws.Range("M5").Value = 0.1124023333333333
$ws.Range("N5").Value = 0.337207
$ws.Range("O5").Value = 0.03490487583665934
$ws.Range("P5").Value = 0.03490487583665934
$ws.Range("Q5").Value = 2.971573966998
$ws.Range("R5").Value = 26.744165702982
$ws.Range("S5").Value = 0.02349006481025905
$ws.Range("T5").Value = 0.02349006481025905

# Row 6
$ws.Range("I6").Value = 0.6729737392616156
$ws.Range("J6").Value = 0.6729737392616155
$ws.Range("O6").Value = 0.1007904982867776
$ws.Range("P6").Value = 0.1007904982867776
$ws.Range("S6").Value = 0.06782935851409415
$ws.Range("T6").Value = 0.06782935851409415

# Row 7
$ws.Range("I7").Value = 0.6729737392616156
$ws.Range("J7").Value = 0.6729737392616155
$ws.Range("M7").Value = 2.783274666666667
$ws.Range("N7").Value = 8.349824
$ws.Range("O7").Value = 0.8643046258765631
$ws.Range("P7").Value = 0.8643046258765632
$ws.Range("Q7").Value = 73.581270932736
$ws.Range("R7").Value = 662.231438394624
$ws.Range("S7").Value = 0.5816543159372625
$ws.Range("T7").Value = 0.5816543159372624

# Row 8
$ws.Range("G8").Value = 12.55394366666667
$ws.Range("H8").Value = 37.661831
$ws.Range("I8").Value = 0.3195707889299884
$ws.Range("J8").Value = 0.3195707889299884
$ws.Range("M8").Value = 0.1124023333333333
$ws.Range("N8").Value = 0.337207
$ws.Range("O8").Value = 0.03490487583665934
$ws.Range("P8").Value = 0.03490487583665934
$ws.Range("Q8").Value = 1.411092560668556
$ws.Range("R8").Value = 12.699833046017
$ws.Range("S8").Value = 0.01115457870862451
$ws.Range("T8").Value = 0.01115457870862451

# Row 9
$ws.Range("G9").Value = 12.55394366666667
$ws.Range("H9").Value = 37.661831
$ws.Range("I9").Value = 0.3195707889299884
$ws.Range("J9").Value = 0.3195707889299884
$ws.Range("O9").Value = 0.1007904982867776
$ws.Range("P9").Value = 0.1007904982867776
$ws.Range("Q9").Value = 4.074637680537889
$ws.Range("R9").Value = 36.671739124841
$ws.Range("S9").Value = 0.03220969905415214
$ws.Range("T9").Value = 0.03220969905415215

# Row 10
$ws.Range("G10").Value = 12.55394366666667
$ws.Range("H10").Value = 37.661831
$ws.Range("I10").Value = 0.3195707889299884
$ws.Range("J10").Value = 0.3195707889299884
$ws.Range("M10").Value = 2.783274666666667
$ws.Range("N10").Value = 8.349824
$ws.Range("O10").Value = 0.8643046258765631
$ws.Range("P10").Value = 0.8643046258765632
$ws.Range("Q10").Value = 34.94107337419378
$ws.Range("R10").Value = 314.469660367744
$ws.Range("S10").Value = 0.2762065111672117
$ws.Range("T10").Value = 0.2762065111672117
